$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..4) {
    $q = $ws.Cells.Item($row, 17)  # Q
    $r = $ws.Cells.Item($row, 18)  # R
    $q.Value = [Math]::Round([double]$q.Value2)
    $r.Value = [Math]::Round([double]$r.Value2)

    $ws.Cells.Item($row, 26).ClearContents()  # Z
    $ws.Cells.Item($row, 28).ClearContents()  # AB
}
